$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column V (V, W, X become new; old V..AG shift to Y..AJ)
$ws.Columns("V:X").Insert()

# New header cells in row 7 (formatting is inherited from the insert, matching style 13)
$ws.Range("V7").Value2 = "KL thùng & gang lỏng (T)"
$ws.Range("W7").Value2 = "KL thùng (T)"
$ws.Range("X7").Value2 = "KL gang lỏng (T)"

# Adjust column widths to match the post-edit layout
$ws.Columns("V").ColumnWidth = 11.166666666666666
$ws.Columns("W").ColumnWidth = 8.833333333333334
$ws.Columns("X").ColumnWidth = 8.666666666666666
$ws.Columns("Y").ColumnWidth = 9.666666666666666
$ws.Columns("AB:AC").ColumnWidth = 9.833333333333334
$ws.Columns("AD:AE").ColumnWidth = 10

# Restore view state (scrolled back to the left, cursor parked at K18)
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("K18").Select()
